$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A (time values like "10.30", "16.00") to be stored as text
# so Excel does not auto-convert them into numbers.
$ws.Range("A1:A10").NumberFormat = "@"

# Header row stays the same: Tid | Bruker | Dato
$ws.Range("A1").Value = "Tid"
$ws.Range("B1").Value = "Bruker"
$ws.Range("C1").Value = "Dato"

# Data rows reflecting the new "Tidsberekning" demo table
$ws.Range("A2").Value = "10.30"
$ws.Range("B2").Value = "Marius Sørenes"
$ws.Range("C2").Value = "27.07.2016"

$ws.Range("A3").Value = "16.00"
$ws.Range("B3").Value = "Johannes  Steinsbø"
$ws.Range("C3").Value = "27.07.2016"

$ws.Range("A4").Value = "18.00"
$ws.Range("B4").Value = "Marius Sørenes"
$ws.Range("C4").Value = "27.07.2016"

$ws.Range("A5").Value = "10.00"
$ws.Range("B5").Value = "Bergliot  Olavsen"
$ws.Range("C5").Value = "28.07.2016"

$ws.Range("A6").Value = "11.00"
$ws.Range("B6").Value = "Johannes  Steinsbø"
$ws.Range("C6").Value = "29.07.2016"

$ws.Range("A7").Value = "12.00"
$ws.Range("B7").Value = "Jon  Olav"
$ws.Range("C7").Value = "29.07.2016"

$ws.Range("A8").Value = "13.00"
$ws.Range("B8").Value = "Nessen "
$ws.Range("C8").Value = "29.07.2016"

$ws.Range("A9").Value = "14.00"
$ws.Range("B9").Value = "Sigmund  Steinsbø"
$ws.Range("C9").Value = "30.07.2016"

$ws.Range("A10").Value = "15.00"
$ws.Range("B10").Value = "Ola  Nordmann"
$ws.Range("C10").Value = "30.07.2016"
